$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# =========================================================================
# Sheet1-5Columns: column D (previously text "...Col4") becomes a date
# column; columns A, B, C, E drop their explicit custom widths (back to
# the sheet default) while D gets its own custom width.
# =========================================================================

$ws1.Range("A1:C1").EntireColumn.ClearFormats()
$ws1.Range("E1").EntireColumn.ClearFormats()

$ws1.Range("D2:D4").NumberFormat = "mm/dd/yy"
$ws1.Range("D2").Value = 45777.5397278588
$ws1.Range("D3").Value = 45777.5397279977
$ws1.Range("D4").Value = 45777.5397280324

$ws1.Range("D1").EntireColumn.ColumnWidth = 9.14

# =========================================================================
# Sheet2-3Columns: drop the two phantom trailing rows, add a new data row
# (Row3), and drop the explicit custom column widths (back to default).
# =========================================================================

$ws2.Range("A1048575:A1048576").EntireRow.Delete()

$ws2.Cells.Item(4, 1).Value = "Sheet2-Row3-ColA"
$ws2.Cells.Item(4, 2).Value = "Sheet2-Row3-ColB"
$ws2.Cells.Item(4, 3).Value = "Sheet2-Data3"

$ws2.Range("A1:C1").EntireColumn.ClearFormats()

# =========================================================================
# Selection / active-sheet bookkeeping (mirrors the target sheetViews).
# =========================================================================

$ws2.Range("A1").Select()
$ws1.Range("F3").Select()
$ws1.Activate()
